# Apply the sushi-config.yaml "pin-canonicals" regeneration changes to the
# EyeColor StructureDefinition workbook:
#   1. Bump the "Date" metadata value to the new publish timestamp.
#   2. Pin the Base Definition canonical URL to its version (|4.0.1).
#   3. Pin the Binding Value Set canonical URL to its version (|0.1.0).
#   4. The "Binding Value Set" column on the Elements sheet grows wider to
#      fit the longer, now-versioned URL (bestFit recalculation).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Date: 2025-12-03T08:52:56+00:00 -> 2025-12-19T08:24:59+00:00
$meta.Range("B8").Value = "2025-12-19T08:24:59+00:00"

# Base Definition: append the pinned FHIR version
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet -----------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set: append the pinned IG version
$elements.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/[code]/ValueSet/EyeColorVS|0.1.0"

# The "Binding Value Set" column (Z) best-fits to the new, longer text -
# widen it to match (character-width grid makes 56.41015625 unreachable
# bit-exactly via COM's pixel-quantized ColumnWidth; land in the closest
# attainable bucket just above the old 51.87890625 width).
$elements.Columns.Item(26).ColumnWidth = 55.5
